$d = $word.ActiveDocument

# "Programa" paragraph (Portuguese): insert a manual line break after each
# of the first two sentences, splitting the single run into three.
$r1 = $d.Content
$ok1 = $r1.Find.Execute("desigualdades. Funções Reais", $true, $false, $false, $false, $false, $true, 1, $false, "desigualdades. ^lFunções Reais", 2)
if (-not $ok1) { throw "Find/Replace 1 (desigualdades./Funções Reais) did not match" }

$r2 = $d.Content
$ok2 = $r2.Find.Execute("hiperbólicas. Modelagem", $true, $false, $false, $false, $false, $true, 1, $false, "hiperbólicas. ^lModelagem", 2)
if (-not $ok2) { throw "Find/Replace 2 (hiperbólicas./Modelagem) did not match" }

# "Programa" paragraph (English, italic): same split, three sentences.
$r3 = $d.Content
$ok3 = $r3.Find.Execute("inequalities.Real Functions", $true, $false, $false, $false, $false, $true, 1, $false, "inequalities.^lReal Functions", 2)
if (-not $ok3) { throw "Find/Replace 3 (inequalities./Real Functions) did not match" }

$r4 = $d.Content
$ok4 = $r4.Find.Execute("functions.Modeling", $true, $false, $false, $false, $false, $true, 1, $false, "functions.^lModeling", 2)
if (-not $ok4) { throw "Find/Replace 4 (functions./Modeling) did not match" }

# "Bibliografia" paragraph: insert a double manual line break between each
# of the four references.
$r5 = $d.Content
$ok5 = $r5.Find.Execute("v.1.ANTON", $true, $false, $false, $false, $false, $true, 1, $false, "v.1.^l^lANTON", 2)
if (-not $ok5) { throw "Find/Replace 5 (v.1./ANTON) did not match" }

$r6 = $d.Content
$ok6 = $r6.Find.Execute("2007.THOMAS", $true, $false, $false, $false, $false, $true, 1, $false, "2007.^l^lTHOMAS", 2)
if (-not $ok6) { throw "Find/Replace 6 (2007./THOMAS) did not match" }

$r7 = $d.Content
$ok7 = $r7.Find.Execute("v.1,FLEMMING", $true, $false, $false, $false, $false, $true, 1, $false, "v.1,^l^lFLEMMING", 2)
if (-not $ok7) { throw "Find/Replace 7 (v.1,/FLEMMING) did not match" }

Write-Output ("All replacements applied: " + $ok1 + "," + $ok2 + "," + $ok3 + "," + $ok4 + "," + $ok5 + "," + $ok6 + "," + $ok7)
